$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 126.22222
$ws.Range("I8").Value = 126.22222
$ws.Range("K8").Value = 378.66666
$ws.Range("M8").Value = -239.66666
$ws.Range("H17").Value = 3948.7542
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 4011.2334
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 12033.7002
$ws.Range("M17").Value = -432
$ws.Range("N17").Value = -12369.7002
$ws.Range("H33").Value = 3684.5454
$ws.Range("I33").Value = 4211.3687
$ws.Range("K33").Value = 4211.3687
$ws.Range("M33").Value = -3982.3687
$ws.Range("H137").Value = 56657.75
$ws.Range("I137").Value = 75027
$ws.Range("J137").Value = 1550
$ws.Range("K137").Value = 225081
$ws.Range("L137").Value = 4650
$ws.Range("M137").Value = -222531
$ws.Range("N137").Value = -9750
$ws.Range("H138").Value = 4596.9473
$ws.Range("I138").Value = 3229.7
$ws.Range("J138").Value = 5085.25
$ws.Range("K138").Value = 9689.099999999999
$ws.Range("L138").Value = 15255.75
$ws.Range("M138").Value = -4549.099999999999
$ws.Range("N138").Value = -25535.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3223.359
$ws.Range("I32").Value = 2157.0747
$ws.Range("K32").Value = 2157.0747
$ws.Range("M32").Value = -1870.0747
$ws.Range("H45").Value = 11069784
$ws.Range("I45").Value = 15985534
$ws.Range("J45").Value = 9346.75
$ws.Range("K45").Value = 15985534
$ws.Range("L45").Value = 9346.75
$ws.Range("M45").Value = -15985157
$ws.Range("N45").Value = -10100.75
$ws.Range("H61").Value = 5543.1562
$ws.Range("I61").Value = 6087.5186
$ws.Range("J61").Value = 2603.6
$ws.Range("K61").Value = 6087.5186
$ws.Range("L61").Value = 2603.6
$ws.Range("M61").Value = -5875.5186
$ws.Range("N61").Value = -3027.6
$ws.Range("H74").Value = 36376.285
$ws.Range("I74").Value = 7788.93
$ws.Range("J74").Value = 108685.47
$ws.Range("K74").Value = 7788.93
$ws.Range("L74").Value = 108685.47
$ws.Range("M74").Value = -6914.93
$ws.Range("N74").Value = -110433.47
$ws.Range("H77").Value = 36376.285
$ws.Range("I77").Value = 7788.93
$ws.Range("J77").Value = 108685.47
$ws.Range("K77").Value = 38944.65
$ws.Range("L77").Value = 543427.35
$ws.Range("M77").Value = -34576.65
$ws.Range("N77").Value = -552163.35
$ws.Range("H102").Value = 2691552
$ws.Range("I102").Value = 3089621
$ws.Range("J102").Value = 4587.5
$ws.Range("K102").Value = 3089621
$ws.Range("L102").Value = 4587.5
$ws.Range("M102").Value = -3087999
$ws.Range("N102").Value = -7831.5
$ws.Range("H136").Value = 5543.1562
$ws.Range("I136").Value = 6087.5186
$ws.Range("J136").Value = 2603.6
$ws.Range("K136").Value = 18262.5558
$ws.Range("L136").Value = 7810.799999999999
$ws.Range("M136").Value = -15712.5558
$ws.Range("N136").Value = -12910.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4356123.5
$ws.Range("I86").Value = 6260699.5
$ws.Range("K86").Value = 6260699.5
$ws.Range("M86").Value = -6259576.5
$ws.Range("H89").Value = 4356123.5
$ws.Range("I89").Value = 6260699.5
$ws.Range("K89").Value = 31303497.5
$ws.Range("M89").Value = -31297881.5
$ws.Range("H96").Value = 19983.777
$ws.Range("I96").Value = 19983.777
$ws.Range("K96").Value = 19983.777
$ws.Range("M96").Value = -17237.777

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 100299.5
$ws.Range("J17").Value = 599
$ws.Range("L17").Value = 599
$ws.Range("N17").Value = -947
$ws.Range("H31").Value = 3624.11
$ws.Range("J31").Value = 4276.924
$ws.Range("L31").Value = 4276.924
$ws.Range("N31").Value = -4866.924
$ws.Range("H34").Value = 3624.11
$ws.Range("J34").Value = 4276.924
$ws.Range("L34").Value = 4276.924
$ws.Range("N34").Value = -4680.924
$ws.Range("H50").Value = 6699.143
$ws.Range("J50").Value = 6699.143
$ws.Range("L50").Value = 6699.143
$ws.Range("N50").Value = -7949.143
$ws.Range("H62").Value = 3399
$ws.Range("J62").Value = 3498.75
$ws.Range("L62").Value = 3498.75
$ws.Range("N62").Value = -4746.75
$ws.Range("H65").Value = 3399
$ws.Range("J65").Value = 3498.75
$ws.Range("L65").Value = 17493.75
$ws.Range("N65").Value = -23733.75
$ws.Range("H105").Value = 4300
$ws.Range("I105").Value = 4300
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 4300
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -2553
$ws.Range("N105").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3548
$ws.Range("I115").Value = 2580
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 7740
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -6565
$ws.Range("N115").Value = -17350
$ws.Range("H122").Value = 746.7143
$ws.Range("I122").Value = 880.6
$ws.Range("J122").Value = 704.875
$ws.Range("K122").Value = 7925.400000000001
$ws.Range("L122").Value = 6343.875
$ws.Range("M122").Value = -5475.400000000001
$ws.Range("N122").Value = -11243.875
$ws.Range("H137").Value = 4641.6875
$ws.Range("I137").Value = 4748.75
$ws.Range("J137").Value = 4606
$ws.Range("K137").Value = 14246.25
$ws.Range("L137").Value = 13818
$ws.Range("M137").Value = -9146.25
$ws.Range("N137").Value = -24018

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6283.4375
$ws.Range("I2").Value = 48.5
$ws.Range("J2").Value = 16675
$ws.Range("K2").Value = 48.5
$ws.Range("L2").Value = 16675
$ws.Range("M2").Value = 64.5
$ws.Range("N2").Value = -16901
$ws.Range("H9").Value = 729.6667
$ws.Range("I9").Value = 599.5
$ws.Range("J9").Value = 990
$ws.Range("K9").Value = 599.5
$ws.Range("L9").Value = 990
$ws.Range("M9").Value = -429.5
$ws.Range("N9").Value = -1330
$ws.Range("H132").Value = 4876.9316
$ws.Range("I132").Value = 2588.5833
$ws.Range("K132").Value = 7765.749899999999
$ws.Range("M132").Value = -5235.749899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 34900
$ws.Range("I4").Value = 34900
$ws.Range("K4").Value = 34900
$ws.Range("M4").Value = -34787
$ws.Range("H22").Value = 51185.89
$ws.Range("I22").Value = 112270
$ws.Range("J22").Value = 2318.6
$ws.Range("K22").Value = 112270
$ws.Range("L22").Value = 2318.6
$ws.Range("M22").Value = -111975
$ws.Range("N22").Value = -2908.6
$ws.Range("H27").Value = 51185.89
$ws.Range("I27").Value = 112270
$ws.Range("J27").Value = 2318.6
$ws.Range("K27").Value = 112270
$ws.Range("L27").Value = 2318.6
$ws.Range("M27").Value = -112163
$ws.Range("N27").Value = -2532.6
$ws.Range("H28").Value = 34900
$ws.Range("I28").Value = 34900
$ws.Range("K28").Value = 34900
$ws.Range("M28").Value = -34668
$ws.Range("H37").Value = 34900
$ws.Range("I37").Value = 34900
$ws.Range("K37").Value = 34900
$ws.Range("M37").Value = -34793
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256
$ws.Range("H82").Value = 4274454
$ws.Range("I82").Value = 7937371.5
$ws.Range("J82").Value = 1050
$ws.Range("K82").Value = 7937371.5
$ws.Range("L82").Value = 1050
$ws.Range("M82").Value = -7937010.5
$ws.Range("N82").Value = -1772
$ws.Range("H85").Value = 4274454
$ws.Range("I85").Value = 7937371.5
$ws.Range("J85").Value = 1050
$ws.Range("K85").Value = 7937371.5
$ws.Range("L85").Value = 1050
$ws.Range("M85").Value = -7936123.5
$ws.Range("N85").Value = -3546
$ws.Range("H136").Value = 28055.744
$ws.Range("I136").Value = 43382.293
$ws.Range("J136").Value = 3533.2666
$ws.Range("K136").Value = 130146.879
$ws.Range("L136").Value = 10599.7998
$ws.Range("M136").Value = -127596.879
$ws.Range("N136").Value = -15699.7998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10394.625
$ws.Range("J45").Value = 10394.625
$ws.Range("L45").Value = 10394.625
$ws.Range("N45").Value = -11376.625
$ws.Range("H62").Value = 14310.914
$ws.Range("I62").Value = 68783.664
$ws.Range("J62").Value = 9204.093999999999
$ws.Range("K62").Value = 68783.664
$ws.Range("L62").Value = 9204.093999999999
$ws.Range("M62").Value = -68159.664
$ws.Range("N62").Value = -10452.094
$ws.Range("H65").Value = 14310.914
$ws.Range("I65").Value = 68783.664
$ws.Range("J65").Value = 9204.093999999999
$ws.Range("K65").Value = 343918.32
$ws.Range("L65").Value = 46020.46999999999
$ws.Range("M65").Value = -340798.32
$ws.Range("N65").Value = -52260.46999999999
$ws.Range("H81").Value = 33337274
$ws.Range("I81").Value = 55557124
$ws.Range("J81").Value = 7501
$ws.Range("K81").Value = 111114248
$ws.Range("L81").Value = 15002
$ws.Range("M81").Value = -111113187
$ws.Range("N81").Value = -17124
$ws.Range("H84").Value = 33337274
$ws.Range("I84").Value = 55557124
$ws.Range("J84").Value = 7501
$ws.Range("K84").Value = 555571240
$ws.Range("L84").Value = 75010
$ws.Range("M84").Value = -555565936
$ws.Range("N84").Value = -85618
$ws.Range("H126").Value = 3858.8125
$ws.Range("I126").Value = 3582.28
$ws.Range("J126").Value = 4846.4287
$ws.Range("K126").Value = 10746.84
$ws.Range("L126").Value = 14539.2861
$ws.Range("M126").Value = -8276.84
$ws.Range("N126").Value = -19479.2861
$ws.Range("H136").Value = 7655.8057
$ws.Range("I136").Value = 10089.479
$ws.Range("K136").Value = 30268.437
$ws.Range("M136").Value = -27718.437
